# Updates cryptos list cell values (prices / 1h volume %) per the
# upstream scrape refresh, and fixes the ranking order for two
# coin pairs whose rows swapped position (SuiNetwork/WrappedeETH
# and BabyDogeCoin/Optimism).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without Excel
# reinterpreting numeric-looking strings (e.g. "5.00") as numbers,
# which would both change the stored value (drop trailing zeros)
# and change the cell type away from text. Number format is reset
# back to General/Normal style afterwards so no visible formatting
# residue is left on the cell.
function Set-TextCell {
    param($Ws, $Addr, $Val)
    $Ws.Range($Addr).NumberFormat = "@"
    $Ws.Range($Addr).Value = $Val
    $Ws.Range($Addr).Style = "Normal"
}

$ws.Range('D2').Value = '68.801.65'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '2.508.48'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextCell $ws 'D5' '591.33'
$ws.Range('E5').Value = '  +0.60%  '
Set-TextCell $ws 'D6' '174.18'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('E7').Value = '  +0.03%  '
Set-TextCell $ws 'D8' '0.516'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '2.509.51'
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E10').Value = '  +5.03%  '
$ws.Range('E11').Value = '  -1.07%  '
Set-TextCell $ws 'D12' '5.00'
$ws.Range('E12').Value = '  +1.10%  '
Set-TextCell $ws 'D13' '0.334'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').Value = '2.950.70'
$ws.Range('E14').Value = '  +1.15%  '
Set-TextCell $ws 'D15' '25.68'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '68.719.50'
$ws.Range('E16').Value = '  +1.78%  '
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = '2.522.43'
$ws.Range('E18').Value = '  +5.89%  '
Set-TextCell $ws 'D19' '361.30'
$ws.Range('E19').Value = '  +2.47%  '
Set-TextCell $ws 'D20' '7.50'
$ws.Range('E20').Value = '  -0.77%  '
Set-TextCell $ws 'D21' '10.87'
$ws.Range('E21').Value = '  -2.25%  '
Set-TextCell $ws 'D22' '4.00'
$ws.Range('E22').Value = '  -2.29%  '
$ws.Range('E23').Value = '  -0.03%  '
Set-TextCell $ws 'D24' '70.24'
$ws.Range('E24').Value = '  -0.55%  '
Set-TextCell $ws 'D25' '4.13'
$ws.Range('E25').Value = '  -4.37%  '
Set-TextCell $ws 'D26' '8.86'
$ws.Range('E26').Value = '  -3.64%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.640.29'
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('B28').Value = 'SuiNetwork'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextCell $ws 'D28' '1.64'
$ws.Range('E28').Value = '  -7.37%  '
Set-TextCell $ws 'D29' '0.999'
$ws.Range('E29').Value = '  -0.02%  '
Set-TextCell $ws 'D30' '504.91'
$ws.Range('E30').Value = '  -1.02%  '
$ws.Range('D31').Value = '0.0₃0873'
$ws.Range('E31').Value = '  -4.68%  '
Set-TextCell $ws 'D32' '7.69'
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('E34').Value = '  -4.53%  '
$ws.Range('E35').Value = '  +0.08%  '
Set-TextCell $ws 'D36' '162.02'
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('E37').Value = '  -4.74%  '
Set-TextCell $ws 'D38' '18.49'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('E40').Value = '  +0.04%  '
Set-TextCell $ws 'D41' '1.30'
$ws.Range('E41').Value = '  -2.95%  '
Set-TextCell $ws 'D42' '1.69'
$ws.Range('E42').Value = '  -3.14%  '
Set-TextCell $ws 'D43' '4.71'
$ws.Range('E43').Value = '  -3.97%  '
Set-TextCell $ws 'D44' '0.316'
$ws.Range('E44').Value = '  -4.48%  '
Set-TextCell $ws 'D45' '2.30'
$ws.Range('E45').Value = '  -5.29%  '
Set-TextCell $ws 'D46' '149.63'
$ws.Range('E46').Value = '  +2.35%  '
Set-TextCell $ws 'D47' '3.53'
$ws.Range('E47').Value = '  -0.04%  '
Set-TextCell $ws 'D48' '0.509'
$ws.Range('E48').Value = '  -1.57%  '
Set-TextCell $ws 'D49' '0.0735'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
Set-TextCell $ws 'D50' '1.56'
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0245'
$ws.Range('E51').Value = '  -5.10%  '
